$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (date number format + borders on column A) from the last existing row
# down onto the newly added row 53 before writing values into it.
$ws.Range("A52:E52").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)

$arr = New-Object 'object[,]' 52,5
$arr[0,0] = 39400
$arr[0,1] = 2007
$arr[0,2] = -0.7015558851707349
$arr[0,3] = 2008
$arr[0,4] = 1.245556224328537
$arr[1,0] = 39583
$arr[1,1] = 2008
$arr[1,2] = 1.247512323997491
$arr[1,3] = 2009
$arr[1,4] = 2.82953744009995
$arr[2,0] = 39765
$arr[2,1] = 2008
$arr[2,2] = 0.3590181115727287
$arr[2,3] = 2009
$arr[2,4] = 2.957479223435766
$arr[3,0] = 39948
$arr[3,1] = 2009
$arr[3,2] = 0.5917823527752164
$arr[3,3] = 2010
$arr[3,4] = -1.985049937499994
$arr[4,0] = 40130
$arr[4,1] = 2009
$arr[4,2] = -0.01587181126743165
$arr[4,3] = 2010
$arr[4,4] = 3.835144370775678
$arr[5,0] = 40310
$arr[5,1] = 2010
$arr[5,2] = 1.985742476411234
$arr[5,3] = 2011
$arr[5,4] = 7.819356632099961
$arr[6,0] = 40494
$arr[6,1] = 2010
$arr[6,2] = -0.02256889165885845
$arr[6,3] = 2011
$arr[6,4] = -0.6757980944263386
$arr[7,0] = 40676
$arr[7,1] = 2011
$arr[7,2] = -0.7518797681958955
$arr[7,3] = 2012
$arr[7,4] = -1.590425574400001
$arr[8,0] = 40862
$arr[8,1] = 2011
$arr[8,2] = 0.09611428386597787
$arr[8,3] = 2012
$arr[8,4] = -1.14533080410083
$arr[9,0] = 41044
$arr[9,1] = 2012
$arr[9,2] = -0.5765930039052902
$arr[9,3] = 2013
$arr[9,4] = -1.590425574400001
$arr[10,0] = 41228
$arr[10,1] = 2012
$arr[10,2] = -0.1827723404408288
$arr[10,3] = 2013
$arr[10,4] = -0.9712385417624092
$arr[11,0] = 41409
$arr[11,1] = 2013
$arr[11,2] = -0.07642926654481963
$arr[11,3] = 2014
$arr[11,4] = -0.3994003998999851
$arr[12,0] = 41592
$arr[12,1] = 2013
$arr[12,2] = -0.001350220946483294
$arr[12,3] = 2014
$arr[12,4] = 0.6008487920565297
$arr[13,0] = 41774
$arr[13,1] = 2014
$arr[13,2] = 0.9274109147535681
$arr[13,3] = 2015
$arr[13,4] = 2.829537440099972
$arr[14,0] = 41957
$arr[14,1] = 2014
$arr[14,2] = -0.075754880139145
$arr[14,3] = 2015
$arr[14,4] = -1.589980933184099
$arr[15,0] = 42137
$arr[15,1] = 2015
$arr[15,2] = -0.5259734324162268
$arr[15,3] = 2016
$arr[15,4] = -1.194610791900008
$arr[16,0] = 42321
$arr[16,1] = 2015
$arr[16,2] = -0.5761528471665445
$arr[16,3] = 2016
$arr[16,4] = 0.4501721032283079
$arr[17,0] = 42503
$arr[17,1] = 2016
$arr[17,2] = 0.07456754038981384
$arr[17,3] = 2017
$arr[17,4] = 0.4006004000999708
$arr[18,0] = 42689
$arr[18,1] = 2016
$arr[18,2] = -0.2011999787958185
$arr[18,3] = 2017
$arr[18,4] = -0.1501751374934579
$arr[19,0] = 42867
$arr[19,1] = 2017
$arr[19,2] = -0.07666472728168339
$arr[19,3] = 2018
$arr[19,4] = -1.590425574400012
$arr[20,0] = 43053
$arr[20,1] = 2017
$arr[20,2] = 0.1213692818849532
$arr[20,3] = 2018
$arr[20,4] = 1.204188598110267
$arr[21,0] = 43145
$arr[21,1] = 2018
$arr[21,2] = 0.1986438914956645
$arr[21,3] = 2019
$arr[21,4] = 0
$arr[22,0] = 43235
$arr[22,1] = 2018
$arr[22,2] = -0.2004754673795017
$arr[22,3] = 2019
$arr[22,4] = -0.3994003998999962
$arr[23,0] = 43326
$arr[23,1] = 2018
$arr[23,2] = 0.1743342151774741
$arr[23,3] = 2019
$arr[23,4] = 1.483403742553846
$arr[24,0] = 43418
$arr[24,1] = 2018
$arr[24,2] = 0.1493219406571766
$arr[24,3] = 2019
$arr[24,4] = 2.372078088364726
$arr[25,0] = 43510
$arr[25,1] = 2019
$arr[25,2] = -0.7530239469328848
$arr[25,3] = 2020
$arr[25,4] = -2.378486270400004
$arr[26,0] = 43600
$arr[26,1] = 2019
$arr[26,2] = -0.7283174404322912
$arr[26,3] = 2020
$arr[26,4] = -2.378486270400004
$arr[27,0] = 43691
$arr[27,1] = 2019
$arr[27,2] = -0.2044553505917923
$arr[27,3] = 2020
$arr[27,4] = 0.8757395670492052
$arr[28,0] = 43783
$arr[28,1] = 2019
$arr[28,2] = -0.4278219446121612
$arr[28,3] = 2020
$arr[28,4] = -2.378564786744741
$arr[29,0] = 43875
$arr[29,1] = 2020
$arr[29,2] = 0.6687116611293176
$arr[29,3] = 2021
$arr[29,4] = 2.42168652960002
$arr[30,0] = 43966
$arr[30,1] = 2020
$arr[30,2] = -0.03096525636255842
$arr[30,3] = 2021
$arr[30,4] = 1.205410808099971
$arr[31,0] = 44068
$arr[31,1] = 2020
$arr[31,2] = -0.4412356890029168
$arr[31,3] = 2021
$arr[31,4] = -1.287161643752965
$arr[32,0] = 44159
$arr[32,1] = 2020
$arr[32,2] = -1.026566979837418
$arr[32,3] = 2021
$arr[32,4] = -3.329887238705409
$arr[33,0] = 44251
$arr[33,1] = 2021
$arr[33,2] = -0.9187270897784594
$arr[33,3] = 2022
$arr[33,4] = -2.100316115087164
$arr[34,0] = 44341
$arr[34,1] = 2021
$arr[34,2] = -0.04074803603358879
$arr[34,3] = 2022
$arr[34,4] = -0.2354831990173722
$arr[35,0] = 44432
$arr[35,1] = 2021
$arr[35,2] = -0.04976849661378902
$arr[35,3] = 2022
$arr[35,4] = -6.367504067266704
$arr[36,0] = 44525
$arr[36,1] = 2021
$arr[36,2] = 0.3179894933462268
$arr[36,3] = 2022
$arr[36,4] = 0.09136665742732752
$arr[37,0] = 44617
$arr[37,1] = 2022
$arr[37,2] = 0.7185676380534911
$arr[37,3] = 2023
$arr[37,4] = 1.072043735370976
$arr[38,0] = 44706
$arr[38,1] = 2022
$arr[38,2] = 1.296559640836992
$arr[38,3] = 2023
$arr[38,4] = 0.439921111559638
$arr[39,0] = 44798
$arr[39,1] = 2022
$arr[39,2] = 0.4784173072842179
$arr[39,3] = 2023
$arr[39,4] = -0.07349332085010429
$arr[40,0] = 44890
$arr[40,1] = 2022
$arr[40,2] = 0.463604920919658
$arr[40,3] = 2023
$arr[40,4] = -1.086119702215593
$arr[41,0] = 44981
$arr[41,1] = 2023
$arr[41,2] = 0.3516957773510798
$arr[41,3] = 2024
$arr[41,4] = -0.7234699379142895
$arr[42,0] = 45071
$arr[42,1] = 2023
$arr[42,2] = 0.8143067496459322
$arr[42,3] = 2024
$arr[42,4] = -0.05551881767973388
$arr[43,0] = 45163
$arr[43,1] = 2023
$arr[43,2] = 0.4279153732809959
$arr[43,3] = 2024
$arr[43,4] = 0.09327584436471525
$arr[44,0] = 45254
$arr[44,1] = 2023
$arr[44,2] = 0.621639092134818
$arr[44,3] = 2024
$arr[44,4] = -0.185520708596465
$arr[45,0] = 45345
$arr[45,1] = 2024
$arr[45,2] = -0.4267958602204081
$arr[45,3] = 2025
$arr[45,4] = 0.00005213214835375624
$arr[46,0] = 45436
$arr[46,1] = 2024
$arr[46,2] = -1.077932696718564
$arr[46,3] = 2025
$arr[46,4] = -0.02296885644381685
$arr[47,0] = 45534
$arr[47,1] = 2024
$arr[47,2] = -0.9378224616154895
$arr[47,3] = 2025
$arr[47,4] = 0.3335355268917262
$arr[48,0] = 45618
$arr[48,1] = 2024
$arr[48,2] = -0.6768900623516982
$arr[48,3] = 2025
$arr[48,4] = 0.9196327481985289
$arr[49,0] = 45713
$arr[49,1] = 2025
$arr[49,2] = 2.173733407444822
$arr[49,3] = 2026
$arr[49,4] = 1.89918335857393
$arr[50,0] = 45800
$arr[50,1] = 2025
$arr[50,2] = 0.4743170952486997
$arr[50,3] = 2026
$arr[50,4] = -0.324019587165425
$arr[51,0] = 45891
$arr[51,1] = 2025
$arr[51,2] = 1.136769786738334
$arr[51,3] = 2026
$arr[51,4] = 0.2972123973886909

$ws.Range("A2:E53").Value = $arr

$ws.Range("A1").Select()
